# Update "Förändrad" (column C) date values from 2023-09-11 (45180) to
# 2023-09-12 (45181) for every data row (rows 2 through 253).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 253

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
